# "Modifique una celda del archivo"
#
# Hoja3!E2 (the "w" input) changes from 2 to 3. Every downstream formula in
# Hoja3 (D7:D10, E7:E10, G7:G10) recalculates automatically from that single
# edit. A new, empty sheet "Hoja5" is appended after Hoja4. The active tab
# moves from Hoja4 to Hoja3 (with E2 selected there), while Hoja4's own
# remembered selection becomes C14.

$wb = $excel.ActiveWorkbook

# 1) The actual data edit: Hoja3!E2 2 -> 3 (drives the formula ripple).
$ws3 = $wb.Worksheets.Item("Hoja3")
$ws3.Range("E2").Value = 3

# 2) Append a new, empty worksheet "Hoja5" as the last tab.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "Hoja5"

# 3) Hoja4 is no longer the active sheet, but it keeps C14 as its selection.
$ws4 = $wb.Worksheets.Item("Hoja4")
$ws4.Range("C14").Select() | Out-Null

# 4) Hoja3 becomes the active/selected tab, with E2 as its selection.
$ws3.Activate()
$ws3.Range("E2").Select() | Out-Null
